$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # TAB1_prevalence_cfr15year_BR
$ws2 = $wb.Worksheets.Item(2)   # TAB2_prevalence_cfr15year_CO
$ws3 = $wb.Worksheets.Item(3)   # TAB3_prevalence_cfr15year_NNE
$ws4 = $wb.Worksheets.Item(4)   # TAB4_prevalence_cfr15year_SSE
$ws5 = $wb.Worksheets.Item(5)   # TAB5_decomposition

# ---------------------------------------------------------------------
# 1. Rename the "Central-West" region label to "Midwest" everywhere it
#    appears (TAB2 column A, rows 4-19, and TAB5 cell C2 which shares
#    the same string).
# ---------------------------------------------------------------------
for ($r = 4; $r -le 19; $r++) {
    $ws2.Cells.Item($r, 1).Value = "Midwest"
}
$ws5.Range("C2").Value = "Midwest"

# ---------------------------------------------------------------------
# 2. TAB2: fix the left-hand region-label border so rows 5-19 match the
#    box formatting already used on row 4 (copy format down).
# ---------------------------------------------------------------------
$ws2.Range("A4").Copy()
$ws2.Range("A5:A19").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# 3. Update remembered selections per sheet and switch the active tab
#    to TAB5_decomposition (last sheet activated/selected stays active).
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("F25").Select()

$ws2.Activate()
$ws2.Range("A1").Select()

$ws5.Activate()
$ws5.Range("C3").Select()
